# Apply the "change index scale for figures" edit:
#  - rescale column A (the income-pooling index) from a 0-1 fraction to a
#    0-100 scale (value * 100)
#  - within each same-index pair of rows, reorder so the "Cohab" row comes
#    before the "Married" row (i.e. swap the non-index columns between the
#    two rows of each pair)
#  - refresh the sheet's sortState to reflect the data range (A2:F59) being
#    sorted on column A
#  - rename the worksheet tab to "fig3"
#  - update the active selection to K17

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab.
$ws.Name = "fig3"

# For every pair of data rows (2,3), (4,5), ... (58,59): the two rows share
# the same index value in column A. Swap columns B:F between the two rows
# of the pair, then rescale column A by *100 for both rows.
for ($r = 2; $r -le 58; $r += 2) {
    $r2 = $r + 1

    $b1 = $ws.Cells.Item($r, 2).Value2
    $c1 = $ws.Cells.Item($r, 3).Value2
    $d1 = $ws.Cells.Item($r, 4).Value2
    $e1 = $ws.Cells.Item($r, 5).Value2
    $f1 = $ws.Cells.Item($r, 6).Value2

    $b2 = $ws.Cells.Item($r2, 2).Value2
    $c2 = $ws.Cells.Item($r2, 3).Value2
    $d2 = $ws.Cells.Item($r2, 4).Value2
    $e2 = $ws.Cells.Item($r2, 5).Value2
    $f2 = $ws.Cells.Item($r2, 6).Value2

    $ws.Cells.Item($r, 2).Value = $b2
    $ws.Cells.Item($r, 3).Value = $c2
    $ws.Cells.Item($r, 4).Value = $d2
    $ws.Cells.Item($r, 5).Value = $e2
    $ws.Cells.Item($r, 6).Value = $f2

    $ws.Cells.Item($r2, 2).Value = $b1
    $ws.Cells.Item($r2, 3).Value = $c1
    $ws.Cells.Item($r2, 4).Value = $d1
    $ws.Cells.Item($r2, 5).Value = $e1
    $ws.Cells.Item($r2, 6).Value = $f1

    $a = $ws.Cells.Item($r, 1).Value2
    $aScaled = $a * 100
    $ws.Cells.Item($r, 1).Value = $aScaled
    $ws.Cells.Item($r2, 1).Value = $aScaled
}

# Re-apply the sort (stable, so row order above is preserved) purely so the
# worksheet's recorded sortState matches the data range A2:F59 / A2:A59.
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("A2:A59"))
$sortObj.SetRange($ws.Range("A1:F59"))
$sortObj.Header = 1
$sortObj.Apply()

# Update the selected cell shown when the workbook is reopened.
$ws.Range("K17").Select()
